# Update last week's in-person attendance count for the "Sept. 18th" seminar
# (Lindsey Novak) on the "2024 - Fall" sheet: 7 -> 6.
# All downstream totals (row I/K columns, the summary rows 17-19, the Q
# column on "Attendance Descriptives", and the chart series) are driven by
# formulas and recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024 - Fall")
$ws.Range("G4").Value = 6
